# Updated symbol list on Mon Dec 26 13:45:38 UTC 2022 with GitHub Actions
#
# This script re-applies the "Price" (column D) refresh values, and a couple
# of "Volume(1h)" (column E) label refreshes, exactly as captured in the
# workbook diff. All of the touched cells hold plain text (not numbers), so
# each write is forced to Text format before the value is set - this avoids
# Excel's automatic numeric coercion (e.g. "242.80" -> 242.8) while writing
# the digits unchanged. The cell style is restored to Normal immediately
# afterwards so no other formatting is affected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Address,
        [string]$Text
    )
    $rng = $ws.Range($Address)
    $rng.NumberFormat = "@"
    $rng.Value = $Text
    $rng.Style = "Normal"
}

# Column D ("Price") updates
Set-TextValue "D2"  "242.80"
Set-TextValue "D4"  "5.389"
Set-TextValue "D5"  "0.05906"
Set-TextValue "D6"  "3.456"
Set-TextValue "D7"  "6.565"
Set-TextValue "D8"  "0.8141"
Set-TextValue "D9"  "0.9179"
Set-TextValue "D10" "0.1415"
Set-TextValue "D11" "0.07430"
Set-TextValue "D12" "0.03251"
Set-TextValue "D13" "0.03047"
Set-TextValue "D14" "0.09345"
Set-TextValue "D15" "3.853"
Set-TextValue "D16" "0.001579"
Set-TextValue "D17" "0.04672"
Set-TextValue "D18" "0.0005953"
Set-TextValue "D19" "0.005889"
Set-TextValue "D20" "0.001291"
Set-TextValue "D21" "0.004901"
Set-TextValue "D22" "0.00009504"
Set-TextValue "D25" "0.3201"
Set-TextValue "D27" "0.0002285"
Set-TextValue "D40" "0.03955"
Set-TextValue "D41" "0.006182"
Set-TextValue "D42" "0.1072"
Set-TextValue "D43" "0.002820"
Set-TextValue "D44" "0.008112"
Set-TextValue "D45" "0.00005208"
Set-TextValue "D47" "0.7904"
Set-TextValue "D49" "0.00002101"
Set-TextValue "D50" "0.0002001"

# Column E ("Volume(1h)") label updates
$ws.Range("E18").Value = "17OneONE"
$ws.Range("E20").Value = "19BitKanKAN"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINWorstin24h"
